# Rever_DailyTrack_BALRAJ_2022.xlsx - "Add files via upload" edit
#
# Fills in 4 new progress-log rows (13th day entry) on the "JAN-22" sheet:
#   Row 24: No=13, Date=24-Jan-2022, Application=RPA RLOGIC,
#           Task="1. P&L has been triggered ...", %=100%, Status=Completed
#   Row 25: Task="2. RPA Management template has been triggered ...", %=100%, Status=Completed
#   Row 26: Task="3. Getting Last month closing balance ...", %=100%, Status=Completed
#   Row 27: Task="4. Some more filtering to be done ...", %=20%, Status=WIP
#
# Number-format constant used below:
#   -4122 == xlPasteFormats (paste only the cell format, not the value)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 24 (new "No" / date / application entry) ----
$ws.Range("A24").Value = 13
$ws.Range("B24").Value = 44585
$ws.Range("B2").Copy()
$ws.Range("B24").PasteSpecial(-4122)
$ws.Range("C24").Value = "RPA RLOGIC"
$ws.Range("D24").Value = "1. P&L has been triggered for the all 3 centers and it is done today"
$ws.Range("E24").Value = 1
$ws.Range("E2").Copy()
$ws.Range("E24").PasteSpecial(-4122)
$ws.Range("F24").Value = "Completed"

# ---- Row 25 ----
$ws.Range("D25").Value = "2. RPA Management template has been triggered for the all 3 centers  and it is done today"
$ws.Range("E25").Value = 1
$ws.Range("E3").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("F25").Value = "Completed"

# ---- Row 26 ----
$ws.Range("D26").Value = "3. Getting Last month closing balance and taking as current month opening balance, this steps has been implemented at GL and etc"
$ws.Range("E26").Value = 1
$ws.Range("E3").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("F26").Value = "Completed"

# ---- Row 27 ----
$ws.Range("D27").Value = "4. Some more filtering to be done at Vendor Ledger file for GL to implement is work in progress"
$ws.Range("E27").Value = 0.2
$ws.Range("E3").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("F27").Value = "WIP"

# ---- Restore the selection left behind by copy/paste, matching the
#      author's final cursor position (D34, scrolled so row 10 is on top) ----
$excel.CutCopyMode = $false
$ws.Range("D34").Select() | Out-Null
